$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $value) {
    $r = $ws.Range($cellRef)
    $r.NumberFormat = "@"
    $r.Value = $value
    $r.Style = "Normal"
}

Set-TextValue 'D2' '26.634.61'
Set-TextValue 'E2' '  -2.44%  '
Set-TextValue 'D3' '1.788.41'
Set-TextValue 'E3' '  -2.10%  '
Set-TextValue 'E4' '  +0.23%  '
Set-TextValue 'D5' '307.86'
Set-TextValue 'E5' '  -1.83%  '
Set-TextValue 'E6' '  +0.21%  '
Set-TextValue 'D7' '0.4551'
Set-TextValue 'E7' '  +1.50%  '
Set-TextValue 'D8' '0.3693'
Set-TextValue 'E8' '  -2.63%  '
Set-TextValue 'D9' '0.07204'
Set-TextValue 'E9' '  -4.17%  '
Set-TextValue 'D10' '0.8528'
Set-TextValue 'E10' '  -3.55%  '
Set-TextValue 'D11' '20.31'
Set-TextValue 'E11' '  -3.38%  '
Set-TextValue 'D12' '1.793.47'
Set-TextValue 'E12' '  -1.11%  '
Set-TextValue 'D13' '5.289'
Set-TextValue 'E13' '  -2.15%  '
Set-TextValue 'D14' '0.07035'
Set-TextValue 'E14' '  -1.37%  '
Set-TextValue 'D15' '6.468'
Set-TextValue 'E15' '  -4.50%  '
Set-TextValue 'D16' '89.96'
Set-TextValue 'E16' '  -5.07%  '
Set-TextValue 'D17' '1.003'
Set-TextValue 'E17' '  +0.36%  '
Set-TextValue 'D18' '0.000008598'
Set-TextValue 'E18' '  -2.32%  '
Set-TextValue 'E19' '  +0.18%  '
Set-TextValue 'E20' '  -4.07%  '
Set-TextValue 'D21' '26.644.41'
Set-TextValue 'E21' '  -2.51%  '
Set-TextValue 'D22' '5.278'
Set-TextValue 'E22' '  +0.38%  '
Set-TextValue 'E23' '  -3.95%  '
Set-TextValue 'D24' '2.008.95'
Set-TextValue 'D25' '1.908'
Set-TextValue 'E25' '  -4.13%  '
Set-TextValue 'D26' '149.70'
Set-TextValue 'E26' '  -1.61%  '
Set-TextValue 'D27' '18.08'
Set-TextValue 'E27' '  -3.15%  '
Set-TextValue 'D28' '2.118'
Set-TextValue 'E28' '  -13.83%  '
Set-TextValue 'D29' '5.179'
Set-TextValue 'E29' '  -3.74%  '
Set-TextValue 'D30' '113.87'
Set-TextValue 'E30' '  -3.90%  '
Set-TextValue 'D31' '0.08803'
Set-TextValue 'E31' '  -0.64%  '
Set-TextValue 'D32' '0.7535'
Set-TextValue 'E32' '  -2.53%  '
Set-TextValue 'D33' '1.150'
Set-TextValue 'E33' '  -3.46%  '
Set-TextValue 'D34' '4.425'
Set-TextValue 'E34' '  -3.49%  '
Set-TextValue 'D35' '2.887'
Set-TextValue 'E35' '  +0.09%  '
Set-TextValue 'E36' '  +0.24%  '
Set-TextValue 'D37' '1.108'
Set-TextValue 'E37' '  -0.32%  '
Set-TextValue 'D38' '0.01932'
Set-TextValue 'E38' '  -3.18%  '
Set-TextValue 'D39' '0.05188'
Set-TextValue 'E39' '  -2.42%  '
Set-TextValue 'D40' '2.885'
Set-TextValue 'E40' '  +0.91%  '
Set-TextValue 'D41' '7.108'
Set-TextValue 'E41' '  -4.30%  '
Set-TextValue 'D42' '2.334'
Set-TextValue 'E42' '  +3.41%  '
Set-TextValue 'D43' '0.5195'
Set-TextValue 'D44' '0.1638'
Set-TextValue 'E44' '  -5.39%  '
Set-TextValue 'D45' '8.427'
Set-TextValue 'E45' '  -4.27%  '
Set-TextValue 'D46' '0.4942'
Set-TextValue 'E46' '  -3.20%  '
Set-TextValue 'D47' '10.23'
Set-TextValue 'E47' '  -5.38%  '
Set-TextValue 'B48' 'Quant'
Set-TextValue 'C48' 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
Set-TextValue 'D48' '104.01'
Set-TextValue 'E48' '  -2.50%  '
Set-TextValue 'B49' 'PaxDollar'
Set-TextValue 'C49' 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
Set-TextValue 'D49' '1.000'
Set-TextValue 'E49' '  +0.22%  '
Set-TextValue 'E50' '  -3.87%  '
Set-TextValue 'D51' '0.06271'
Set-TextValue 'E51' '  -1.74%  '
